$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 443.9
$ws.Range("I8").Value = 48.77778
$ws.Range("K8").Value = 146.33334
$ws.Range("M8").Value = -7.333339999999993
$ws.Range("H17").Value = 435.0566
$ws.Range("J17").Value = 245.73334
$ws.Range("L17").Value = 737.20002
$ws.Range("N17").Value = -1073.20002
$ws.Range("H132").Value = 34832840
$ws.Range("I132").Value = 37042124
$ws.Range("K132").Value = 111126372
$ws.Range("M132").Value = -111123842
$ws.Range("H137").Value = 6807064.5
$ws.Range("I137").Value = 11905662
$ws.Range("J137").Value = 8933.333000000001
$ws.Range("K137").Value = 35716986
$ws.Range("L137").Value = 26799.999
$ws.Range("M137").Value = -35714436
$ws.Range("N137").Value = -31899.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1750
$ws.Range("I61").Value = 1750
$ws.Range("K61").Value = 1750
$ws.Range("M61").Value = -1538
$ws.Range("H74").Value = 2711.0667
$ws.Range("I74").Value = 1653.2727
$ws.Range("K74").Value = 1653.2727
$ws.Range("M74").Value = -779.2727
$ws.Range("H77").Value = 2711.0667
$ws.Range("I77").Value = 1653.2727
$ws.Range("K77").Value = 8266.363499999999
$ws.Range("M77").Value = -3898.363499999999
$ws.Range("H132").Value = 7049.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7049.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 21148.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -26208.5
$ws.Range("H136").Value = 1750
$ws.Range("I136").Value = 1750
$ws.Range("K136").Value = 5250
$ws.Range("M136").Value = -2700
$ws.Range("H139").Value = 41347.97
$ws.Range("J139").Value = 41347.97
$ws.Range("L139").Value = 41347.97
$ws.Range("N139").Value = -51627.97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 17600
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 17600
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 17600
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -17944
$ws.Range("H134").Value = 4243.067
$ws.Range("I134").Value = 1863.2
$ws.Range("J134").Value = 9002.799999999999
$ws.Range("K134").Value = 5589.6
$ws.Range("L134").Value = 27008.4
$ws.Range("M134").Value = -3054.6
$ws.Range("N134").Value = -32078.4
$ws.Range("H138").Value = 40796.297
$ws.Range("J138").Value = 40796.297
$ws.Range("L138").Value = 40796.297
$ws.Range("N138").Value = -51076.297
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 16872
$ws.Range("I12").Value = 1612.5
$ws.Range("J12").Value = 24501.75
$ws.Range("K12").Value = 1612.5
$ws.Range("L12").Value = 24501.75
$ws.Range("M12").Value = -1442.5
$ws.Range("N12").Value = -24841.75
$ws.Range("H48").Value = 37960
$ws.Range("J48").Value = 37960
$ws.Range("L48").Value = 37960
$ws.Range("N48").Value = -38912
$ws.Range("H103").Value = 16121.777
$ws.Range("I103").Value = 6585.143
$ws.Range("K103").Value = 6585.143
$ws.Range("M103").Value = -5413.143
$ws.Range("H123").Value = 40887.777
$ws.Range("J123").Value = 40887.777
$ws.Range("L123").Value = 40887.777
$ws.Range("N123").Value = -50687.777
$ws.Range("H138").Value = 45312.234
$ws.Range("J138").Value = 45312.234
$ws.Range("L138").Value = 45312.234
$ws.Range("N138").Value = -55592.234
$ws.Range("H140").Value = 125651.43
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 125651.43
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 125651.43
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -136011.43
$ws.Range("H141").Value = 27954.166
$ws.Range("J141").Value = 27954.166
$ws.Range("L141").Value = 27954.166
$ws.Range("N141").Value = -38314.166
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 787002.2
$ws.Range("I5").Value = 617.3333
$ws.Range("J5").Value = 1215939.4
$ws.Range("K5").Value = 1851.9999
$ws.Range("L5").Value = 3647818.2
$ws.Range("M5").Value = -1739.9999
$ws.Range("N5").Value = -3648042.2
$ws.Range("H10").Value = 264.2
$ws.Range("I10").Value = 264.2
$ws.Range("K10").Value = 792.5999999999999
$ws.Range("M10").Value = -653.5999999999999
$ws.Range("H113").Value = 4464842.5
$ws.Range("I113").Value = 580.86664
$ws.Range("K113").Value = 1742.59992
$ws.Range("M113").Value = 427.4000800000001
$ws.Range("H129").Value = 3277.9
$ws.Range("J129").Value = 5174.5
$ws.Range("L129").Value = 15523.5
$ws.Range("N129").Value = -25523.5
$ws.Range("H135").Value = 787002.2
$ws.Range("I135").Value = 617.3333
$ws.Range("J135").Value = 1215939.4
$ws.Range("K135").Value = 5555.9997
$ws.Range("L135").Value = 10943454.6
$ws.Range("M135").Value = -3020.9997
$ws.Range("N135").Value = -10948524.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 5377.375
$ws.Range("I9").Value = 4505.25
$ws.Range("J9").Value = 6249.5
$ws.Range("K9").Value = 4505.25
$ws.Range("L9").Value = 6249.5
$ws.Range("M9").Value = -4335.25
$ws.Range("N9").Value = -6589.5
$ws.Range("H102").Value = 2114.9395
$ws.Range("I102").Value = 1375.9584
$ws.Range("J102").Value = 4085.5557
$ws.Range("K102").Value = 1375.9584
$ws.Range("L102").Value = 4085.5557
$ws.Range("M102").Value = 246.0416
$ws.Range("N102").Value = -7329.5557
$ws.Range("H132").Value = 6666.3335
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 7599.6
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 22798.8
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -27858.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1402.2727
$ws.Range("I46").Value = 846.4286
$ws.Range("J46").Value = 2375
$ws.Range("K46").Value = 846.4286
$ws.Range("L46").Value = 2375
$ws.Range("M46").Value = -658.4286
$ws.Range("N46").Value = -2751
$ws.Range("H132").Value = 4099.793
$ws.Range("I132").Value = 3134.652
$ws.Range("J132").Value = 7799.5
$ws.Range("K132").Value = 9403.956
$ws.Range("L132").Value = 23398.5
$ws.Range("M132").Value = -6873.956
$ws.Range("N132").Value = -28458.5
$ws.Range("H136").Value = 6027.273
$ws.Range("J136").Value = 8400
$ws.Range("L136").Value = 25200
$ws.Range("N136").Value = -30300
$ws.Range("H139").Value = 48498.332
$ws.Range("J139").Value = 48498.332
$ws.Range("L139").Value = 48498.332
$ws.Range("N139").Value = -58778.332
$ws.Range("H140").Value = 59072.055
$ws.Range("J140").Value = 59072.055
$ws.Range("L140").Value = 59072.055
$ws.Range("N140").Value = -69432.05499999999
$ws.Range("H141").Value = 32125.79
$ws.Range("J141").Value = 32125.79
$ws.Range("L141").Value = 32125.79
$ws.Range("N141").Value = -42485.79
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 28892.309
$ws.Range("J15").Value = 28892.309
$ws.Range("L15").Value = 28892.309
$ws.Range("N15").Value = -29468.309
$ws.Range("H54").Value = 14494.733
$ws.Range("J54").Value = 14494.733
$ws.Range("L54").Value = 14494.733
$ws.Range("N54").Value = -15534.733
$ws.Range("H122").Value = 3737.5405
$ws.Range("I122").Value = 2603.9
$ws.Range("J122").Value = 5071.2354
$ws.Range("K122").Value = 7811.700000000001
$ws.Range("L122").Value = 15213.7062
$ws.Range("M122").Value = -5361.700000000001
$ws.Range("N122").Value = -20113.7062
$ws.Range("H135").Value = 36690.832
$ws.Range("J135").Value = 36690.832
$ws.Range("L135").Value = 36690.832
$ws.Range("N135").Value = -46830.832
$ws.Range("H136").Value = 5758.6
$ws.Range("I136").Value = 2228.9
$ws.Range("J136").Value = 9288.299999999999
$ws.Range("K136").Value = 6686.700000000001
$ws.Range("L136").Value = 27864.9
$ws.Range("M136").Value = -4136.700000000001
$ws.Range("N136").Value = -32964.89999999999
$ws.Range("H138").Value = 41999.332
$ws.Range("J138").Value = 41999.332
$ws.Range("L138").Value = 41999.332
$ws.Range("N138").Value = -52279.332
$ws.Range("H139").Value = 41726.668
$ws.Range("J139").Value = 41726.668
$ws.Range("L139").Value = 41726.668
$ws.Range("N139").Value = -52006.668
$ws.Range("H140").Value = 30343.555
$ws.Range("J140").Value = 30343.555
$ws.Range("L140").Value = 30343.555
$ws.Range("N140").Value = -40703.555
